$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.490.82'
$ws.Range('E2').Value = '  -4.07%  '

$ws.Range('D3').Value = '2.536.38'
$ws.Range('E3').Value = '  -3.96%  '

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '505.69'
$ws.Range('E5').Value = '  -4.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.08'
$ws.Range('E6').Value = '  -7.89%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.16%  '

$ws.Range('E8').Value = '  -4.51%  '

$ws.Range('D9').Value = '2.544.93'
$ws.Range('E9').Value = '  -4.02%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.09'
$ws.Range('E10').Value = '  -8.27%  '

$ws.Range('E11').Value = '  -6.48%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.331'
$ws.Range('E12').Value = '  -5.81%  '

$ws.Range('E13').Value = '  -0.55%  '

$ws.Range('D14').Value = '2.986.13'
$ws.Range('E14').Value = '  -3.64%  '

$ws.Range('D15').Value = '58.485.08'
$ws.Range('E15').Value = '  -4.05%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.66'
$ws.Range('E16').Value = '  -5.63%  '

$ws.Range('E17').Value = '  -5.99%  '

$ws.Range('D18').Value = '2.544.42'
$ws.Range('E18').Value = '  -3.60%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.52'
$ws.Range('E19').Value = '  -4.98%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '338.91'
$ws.Range('E20').Value = '  -4.21%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.08'
$ws.Range('E21').Value = '  -5.41%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.37%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.95'
$ws.Range('E23').Value = '  -4.42%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.64'
$ws.Range('E24').Value = '  -1.52%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.411'
$ws.Range('E25').Value = '  -4.56%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.12%  '

$ws.Range('E27').Value = '  -5.15%  '

$ws.Range('D28').Value = '2.654.92'
$ws.Range('E28').Value = '  -3.57%  '

$ws.Range('D29').Value = '0.0₃0785'
$ws.Range('E29').Value = '  -9.23%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.95'
$ws.Range('E30').Value = '  -6.07%  '

$ws.Range('E31').Value = '  +0.04%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '149.74'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.83'
$ws.Range('E33').Value = '  -5.16%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.50'
$ws.Range('E34').Value = '  -5.07%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('E35').Value = '  -5.78%  '

$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.89'
$ws.Range('E36').Value = '  -6.47%  '

$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.909'
$ws.Range('E37').Value = '  +1.22%  '

$ws.Range('E38').Value = '  -7.72%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.08'
$ws.Range('E39').Value = '  -1.38%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.820'
$ws.Range('E40').Value = '  -11.37%  '

$ws.Range('E41').Value = '  -7.27%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '282.32'
$ws.Range('E42').Value = '  -8.67%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.52'
$ws.Range('E43').Value = '  -7.89%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0996'
$ws.Range('E44').Value = '  -2.32%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  +0.17%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.600'
$ws.Range('E46').Value = '  -6.50%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0531'
$ws.Range('E47').Value = '  -5.65%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.73'
$ws.Range('E48').Value = '  -5.27%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.30'
$ws.Range('E49').Value = '  -0.45%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0226'
$ws.Range('E50').Value = '  -5.53%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.51'
$ws.Range('E51').Value = '  -9.00%  '
